# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" column (D) for the fa783e4a row (row 7)
# on both the zh-cn and de-de status sheets, reflecting a freshly generated
# handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value = "2016-03-08 22:49:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value = "2016-03-08 22:49:17"
